# Update pressure recovery results with newly calculated values
# (new methods added to PressureRecoveryCalculationModel produced refreshed outputs)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Well_1
$ws.Range("B2").Value = 35.47977156947471
$ws.Range("C2").Value = 0.214186018906245
$ws.Range("D2").Value = 3.198007313398188
$ws.Range("E2").Value = 4.908256941102312
$ws.Range("F2").Value = 0.1
$ws.Range("G2").Value = 0.00009963147983989588

# Row 3 - Well_2
$ws.Range("B3").Value = 38.50354809053437
$ws.Range("C3").Value = 0.1103535385144087
$ws.Range("D3").Value = 8.41146270338799
$ws.Range("E3").Value = 3.012861714098561
$ws.Range("F3").Value = 0.1
$ws.Range("G3").Value = 0.00009175455446519026

# Row 4 - Well_3
$ws.Range("B4").Value = 87.00126985761899
$ws.Range("C4").Value = 0.2078844216156235
$ws.Range("D4").Value = 7.571354924040918
$ws.Range("E4").Value = 4.125931442533903
$ws.Range("F4").Value = 0.1
$ws.Range("G4").Value = 0.00008983454907884047

# Row 5 - Well_4
$ws.Range("B5").Value = 83.66670091480998
$ws.Range("C5").Value = 0.1967328393093729
$ws.Range("D5").Value = 4.998921603630805
$ws.Range("E5").Value = 3.694607817391221
$ws.Range("F5").Value = 0.1
$ws.Range("G5").Value = 0.00006929819861501816

# Row 6 - Well_5 (Recovery_Time becomes blank/unavailable)
$ws.Range("B6").Value = 57.55922807493936
$ws.Range("C6").Value = 0.236031542827423
$ws.Range("D6").Value = 2.994633142914468
$ws.Range("E6").Value = -1.818094849856631
$ws.Range("F6").Value = 0.1
$ws.Range("G6").Value = ""

# Row 7 - Well_6
$ws.Range("B7").Value = 13.59258880289455
$ws.Range("C7").Value = 0.1974594140662178
$ws.Range("D7").Value = 8.304780561235734
$ws.Range("E7").Value = 1.628464806417188
$ws.Range("F7").Value = 0.1
$ws.Range("G7").Value = 0.0001661211819431689

# Row 8 - Well_7
$ws.Range("B8").Value = 25.74703298991895
$ws.Range("C8").Value = 0.2275232222525885
$ws.Range("D8").Value = 7.149764342256065
$ws.Range("E8").Value = 3.296464609462556
$ws.Range("F8").Value = 0.1
$ws.Range("G8").Value = 0.0001536996113362175

# Row 9 - Well_8
$ws.Range("B9").Value = 55.18863972267953
$ws.Range("C9").Value = 0.1520796330127055
$ws.Range("D9").Value = 4.467167236381879
$ws.Range("E9").Value = 1.737000730386787
$ws.Range("F9").Value = 0.1
$ws.Range("G9").Value = 0.00005414847891144132

# Row 10 - Well_9
$ws.Range("B10").Value = 67.99750151566482
$ws.Range("C10").Value = 0.1655000430417875
$ws.Range("D10").Value = 3.970055544465525
$ws.Range("E10").Value = -0.1518632039991221
$ws.Range("F10").Value = 0.1
$ws.Range("G10").Value = 0.0000267058151078829

# Row 11 - Well_10
$ws.Range("B11").Value = 18.38067398225872
$ws.Range("C11").Value = 0.1381669083124362
$ws.Range("D11").Value = 8.635658677597409
$ws.Range("E11").Value = -0.8756704885222466
$ws.Range("F11").Value = 0.1
$ws.Range("G11").Value = 0.00002650197059181887
